# Update Name of Algo
# Apply the updated KNN imputation results to the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C11"  = -12.693
    "A12"  = -21.542
    "C23"  = -12.309
    "C28"  = -12.881
    "A32"  = -21.977
    "C32"  = -13.538
    "C34"  = -12.101
    "A36"  = -20.43
    "A38"  = -19.83
    "C42"  = -12.36
    "A46"  = -21.864
    "A54"  = -21.917
    "C54"  = -12.803
    "A55"  = -22.184
    "A67"  = -21.577
    "A69"  = -21.503
    "A72"  = -21.689
    "A91"  = -20.755
    "C97"  = -11.516
    "A99"  = -21.708
    "C99"  = -12.418
    "C101" = -12.188
    "A104" = -21.175
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
